$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.02"
$ws.Range("E2").Value = "'-4.76%"
$ws.Range("D3").Value = "'49.39"
$ws.Range("E3").Value = "'-0.50%"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'-2.65%"
$ws.Range("D5").Value = "'0.07736"
$ws.Range("E5").Value = "'-5.20%"
$ws.Range("D6").Value = "'4.517"
$ws.Range("E6").Value = "'-1.70%"
$ws.Range("E7").Value = "'12.36%"
$ws.Range("D8").Value = "'1.547"
$ws.Range("E8").Value = "'-8.19%"
$ws.Range("D9").Value = "'0.1227"
$ws.Range("E9").Value = "'-9.14%"
$ws.Range("D10").Value = "'0.1958"
$ws.Range("E10").Value = "'-0.63%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04656"
$ws.Range("E11").Value = "'5.40%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09353"
$ws.Range("E12").Value = "'-3.62%"
$ws.Range("D13").Value = "'0.1046"
$ws.Range("E13").Value = "'-0.16%"
$ws.Range("D14").Value = "'0.001261"
$ws.Range("E14").Value = "'-4.97%"
$ws.Range("D15").Value = "'0.04185"
$ws.Range("E15").Value = "'-2.67%"
$ws.Range("D16").Value = "'0.005828"
$ws.Range("E17").Value = "'2,022.74%"
$ws.Range("E18").Value = "'-1.27%"
$ws.Range("D19").Value = "'2.240"
$ws.Range("E19").Value = "'-8.13%"
$ws.Range("D21").Value = "'8.006"
$ws.Range("D22").Value = "'0.1340"
$ws.Range("E22").Value = "'-3.62%"
$ws.Range("E24").Value = "'-2.37%"
$ws.Range("D25").Value = "'0.004012"
$ws.Range("E25").Value = "'-6.00%"
$ws.Range("E26").Value = "'0.32%"
$ws.Range("D38").Value = "'0.02592"
$ws.Range("E38").Value = "'-5.85%"
$ws.Range("D39").Value = "'0.05820"
$ws.Range("E39").Value = "'3.59%"
$ws.Range("D40").Value = "'0.01076"
$ws.Range("E40").Value = "'70.87%"
$ws.Range("D41").Value = "'0.007934"
$ws.Range("E41").Value = "'2.98%"
$ws.Range("D42").Value = "'0.1418"
$ws.Range("E42").Value = "'-2.19%"
$ws.Range("D43").Value = "'0.008460"
$ws.Range("E43").Value = "'10.17%"
$ws.Range("D44").Value = "'0.007711"
$ws.Range("E44").Value = "'-4.82%"
$ws.Range("D45").Value = "'0.3379"
$ws.Range("E45").Value = "'5.80%"
$ws.Range("D46").Value = "'0.00007027"
$ws.Range("E46").Value = "'0.98%"
$ws.Range("E47").Value = "'0.28%"
$ws.Range("D48").Value = "'0.04876"
$ws.Range("E48").Value = "'-20.49%"
$ws.Range("D49").Value = "'0.002626"
$ws.Range("E49").Value = "'-34.32%"
$ws.Range("E50").Value = "'0.28%"
$ws.Range("E51").Value = "'0.28%"
